$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new daily price record was inserted at row 34, shifting all
# subsequent records down by one row (old row 34 -> new row 35, etc.)
$ws.Rows("34:34").Insert()

# Populate the newly inserted row with the new record's data
$ws.Range("A34").Value = 1
$ws.Range("B34").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C34").Value = "Arica y Parinacota"
$ws.Range("D34").Value = 44802
$ws.Range("E34").Value = 15
$ws.Range("F34").Value = "Fruta"
$ws.Range("G34").Value = 100108
$ws.Range("H34").Value = "Tropicales y subtropicales"
$ws.Range("I34").Value = 100108002
$ws.Range("J34").Value = "Mango"
$ws.Range("K34").Value = "Sin especificar"
$ws.Range("L34").Value = "Extra"
$ws.Range("M34").Value = 228
$ws.Range("N34").Value = 9500
$ws.Range("O34").Value = 10000
$ws.Range("P34").Value = 9750
$ws.Range("Q34").Value = "$/bandeja 4 kilos"
$ws.Range("R34").Value = "Brasil"
$ws.Range("S34").Value = 2438
$ws.Range("T34").Value = 4
